$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "3.0.2"
$ws.Range("C2").Value = "3.0.2"
$ws.Range("D2").Value = "3.0.2"

$ws.Range("B31").Value = "Contributor role. At most one PrincipalInvestigator and at least one CorrespondingAuthor are required. Creator cannot be used by itself and requires another role. These roles are provided by the Data Cite schema. Options are:
PrincipalInvestigator
Creator
CoInvestigator
CorrespondingAuthor
ContactPerson
DataCollector
DataCurator
DataManager
Distributor
Editor
HostingInstitution
Producer
ProjectLeader
ProjectManager
ProjectMember
RegistrationAgency
RegistrationAuthority
RelatedPerson
Researcher
ResearchGroup
RightsHolder
Sponsor
Supervisor
WorkPackageLeader
Other."
